$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B,C,D,E,F,I,J,K,L,M,N (column H has no data)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# Each inner array holds the new values for row r, columns in $cols order
$data = @{
    2 = @(1.02, 1.028305176942251, 1.037840886292024, 1.031938115144355, 1.046297517379998, 1.034281192670563, 1.03345823320821, 1.040630621576836, 1.034744788759741, 1.049063347252654, 1.015140285291993)
    3 = @(1.02, 1.029239457544107, 1.038559601824744, 1.032819993698733, 1.047148944768834, 1.034435212746929, 1.034033144457662, 1.04115946888435, 1.035435139295659, 1.049726285302669, 1.015333464076639)
    4 = @(1.02, 1.029844588848975, 1.039024979043206, 1.033391568531935, 1.047700526418051, 1.034533678277589, 1.034405131790107, 1.041501317872235, 1.035882147632458, 1.050155257396035, 1.015458383820889)
    5 = @(1.02, 1.030099126035331, 1.039220698759859, 1.033632082300019, 1.047932565636267, 1.034574786333646, 1.034561509749921, 1.041644946105406, 1.036070141848781, 1.05033559755821, 1.015510880553147)
    6 = @(1.02, 1.030141872126207, 1.039253565323041, 1.033672478696669, 1.047971535043195, 1.034581671718925, 1.034587765943857, 1.041669056895085, 1.036101711087085, 1.050365877443136, 1.015519693828346)
    7 = @(1.02, 1.029847989438622, 1.03902759396633, 1.033394781411453, 1.04770362633319, 1.034534228693313, 1.034407221343032, 1.041503237375473, 1.035884659339314, 1.050157667110641, 1.015459085362045)
    8 = @(1.02, 1.028620798774823, 1.03808371235328, 1.032235954543607, 1.046585125453335, 1.034333491713284, 1.033652530400522, 1.040809419898882, 1.034978032005405, 1.049287388044505, 1.015205587251458)
    9 = @(1.02, 1.026462893412457, 1.036422992106658, 1.030201219588018, 1.044619242258252, 1.033970635351037, 1.032322571362565, 1.039584190646869, 1.033382830497585, 1.047753956964096, 1.01475829787706)
    10 = @(1.02, 1.025027420588707, 1.035317629798499, 1.028849699800504, 1.043312149429177, 1.033722625917546, 1.031435928139813, 1.0387656681741, 1.032321041601073, 1.046731824385925, 1.014459730450078)
    11 = @(1.02, 1.024406600723077, 1.034839437919546, 1.028265673395345, 1.042747012897803, 1.033613794274088, 1.031052013716292, 1.03841084921059, 1.031861688213268, 1.046289283054354, 1.014330363721912)
    12 = @(1.02, 1.02417611417871, 1.034661883374991, 1.028048920267137, 1.042537224309632, 1.033573153282687, 1.030909413006208, 1.038278995678173, 1.031691126407484, 1.046124912153467, 1.014282298781717)
    13 = @(1.02, 1.024225549168068, 1.034699966361821, 1.028095406368093, 1.042582218843401, 1.03358188068208, 1.030940001226822, 1.038307281313747, 1.031727709647036, 1.046160169856887, 1.01429260942819)
    14 = @(1.02, 1.024387546310203, 1.034824759825681, 1.028247752828276, 1.042729669081558, 1.033610439281651, 1.031040226246388, 1.038399951332566, 1.031847588239563, 1.046275695920243, 1.014326390909128)
    15 = @(1.02, 1.024487373179076, 1.034901658148768, 1.02834164244187, 1.0428205350859, 1.033628006562344, 1.031101978510542, 1.038457040741517, 1.031921457713031, 1.046346876526903, 1.014347203166385)
    16 = @(1.02, 1.025068638232885, 1.035349375161015, 1.028888484914164, 1.043349673582176, 1.033729818383351, 1.031461407529007, 1.038789208148904, 1.032351536072617, 1.046761195537795, 1.01446831433595)
    17 = @(1.02, 1.025433451779364, 1.035630334358317, 1.029231824497804, 1.043681815074389, 1.033793296599659, 1.031686870693527, 1.038997463474123, 1.032621422977525, 1.047021100908736, 1.014544261683415)
    18 = @(1.02, 1.025646313580869, 1.035794255195891, 1.029432203447245, 1.043875628846093, 1.033830183207105, 1.031818380241609, 1.039118897255832, 1.032778882751215, 1.047172703936259, 1.014588552252036)
    19 = @(1.02, 1.025718906175313, 1.035850155058797, 1.029500546916376, 1.043941728068235, 1.033842736965701, 1.03186322168689, 1.03916029655324, 1.032832579096814, 1.047224397394703, 1.014603652782457)
    20 = @(1.02, 1.025394303259729, 1.035600185730792, 1.029194975527463, 1.043646171019186, 1.033786500377955, 1.031662680550173, 1.038975123573734, 1.032592462592549, 1.046993215034549, 1.014536114100586)
    21 = @(1.02, 1.024339839045533, 1.034788009393484, 1.028202885604871, 1.042686245107551, 1.033602035446064, 1.031010712407326, 1.03837266391303, 1.031812285247415, 1.046241676124513, 1.014316443443417)
    22 = @(1.02, 1.02367751327933, 1.034277751503809, 1.02758016305281, 1.042083444575439, 1.033484805166439, 1.030600807390997, 1.037993538766269, 1.031322118941362, 1.045769204505013, 1.014178256370424)
    23 = @(1.02, 1.024028562013094, 1.034548211430413, 1.027910180652192, 1.042402929586003, 1.033547069440496, 1.030818104204917, 1.038194551544452, 1.03158193065652, 1.046019665401495, 1.014251518587127)
    24 = @(1.02, 1.025411992567513, 1.035613808467952, 1.029211625636514, 1.043662276782548, 1.033789571728047, 1.031673611032641, 1.038985218130718, 1.032605548422596, 1.047005815446617, 1.014539795668043)
    25 = @(1.02, 1.027020216556612, 1.036852020177368, 1.03072637753895, 1.045126860959887, 1.034065521049629, 1.032666402838615, 1.039901246926823, 1.033794937612152, 1.048150364199649, 1.014874000698027)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}
